# Updates the "Price" (column D) and "Volume(1h)" (column E) columns of the
# cryptos sheet with freshly scraped values, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.186.62'
$ws.Cells.Item(2, 5).Value = '  -1.15%  '
$ws.Cells.Item(3, 4).Value = '1.796.34'
$ws.Cells.Item(3, 5).Value = '  -1.44%  '
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.14%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '314.41'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.28%  '
$ws.Cells.Item(6, 5).Value = '  +0.15%  '
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.5207'
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +2.02%  '
$ws.Cells.Item(8, 5).Value = '  -3.36%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.07978'
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -2.51%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '41.43'
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '1.095'
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -1.26%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '6.296'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.61%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.13%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '20.52'
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -2.77%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '7.298'
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -2.86%  '
$ws.Cells.Item(16, 4).Value = '1.797.15'
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '92.06'
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.54%  '
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '0.00001091'
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -3.66%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '0.06569'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -1.35%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.24%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '17.33'
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -2.68%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '5.952'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -2.32%  '
$ws.Cells.Item(23, 4).Value = '28.218.91'
$ws.Cells.Item(23, 5).Value = '  -1.12%  '
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '11.14'
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -2.09%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '160.40'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +2.47%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '20.47'
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -4.22%  '
$ws.Cells.Item(28, 4).Value = '1.997.19'
$ws.Cells.Item(28, 5).Value = '  -1.64%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '2.341'
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -2.61%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '122.80'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -2.52%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '0.1078'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -1.49%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '1.053'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -5.62%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '3.673'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +0.48%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '5.543'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -3.87%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '0.07245'
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +2.47%  '
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '12.16'
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +7.58%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '0.02321'
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -1.48%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '0.2145'
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -3.74%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '8.720'
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -1.37%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '5.070'
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -3.63%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '0.6159'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -2.46%  '
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '1.163'
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -1.49%  '
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '1.359'
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -2.80%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '13.26'
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -2.15%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '3.770'
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +1.00%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '0.5961'
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.50%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '128.23'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +2.51%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '1.236'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +4.26%  '
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '1.921'
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -3.34%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.06748'
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '73.01'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -1.68%  '
